$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '61.846.41'
$ws.Range("E2").Value = '  +4.53%  '
$ws.Range("D3").Value = '3.073.86'
$ws.Range("E3").Value = '  +2.82%  '
$ws.Range("E4").Value = '  -0.01%  '
$style = $ws.Range("D5").Style
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '579.54'
$ws.Range("D5").Style = $style
$ws.Range("E5").Value = '  +3.16%  '
$style = $ws.Range("D6").Style
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '141.81'
$ws.Range("D6").Style = $style
$ws.Range("E6").Value = '  +3.07%  '
$ws.Range("E7").Value = '  -0.01%  '
$ws.Range("D8").Value = '3.064.79'
$ws.Range("E8").Value = '  +2.87%  '
$ws.Range("E9").Value = '  +1.16%  '
$ws.Range("E10").Value = '  +5.65%  '
$style = $ws.Range("D11").Style
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '5.66'
$ws.Range("D11").Style = $style
$ws.Range("E11").Value = '  +11.11%  '
$style = $ws.Range("D12").Style
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.465'
$ws.Range("D12").Style = $style
$ws.Range("E12").Value = '  +2.21%  '
$style = $ws.Range("D13").Style
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000240'
$ws.Range("D13").Style = $style
$ws.Range("E13").Value = '  +4.55%  '
$style = $ws.Range("D14").Style
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '35.16'
$ws.Range("D14").Style = $style
$ws.Range("E14").Value = '  +4.64%  '
$ws.Range("D16").Value = '3.582.75'
$ws.Range("E17").Value = '  +0.02%  '
$ws.Range("E18").Value = '  +2.64%  '
$ws.Range("D19").Value = '61.765.39'
$ws.Range("E19").Value = '  +4.43%  '
$style = $ws.Range("D20").Style
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '447.10'
$ws.Range("D20").Style = $style
$ws.Range("E20").Value = '  +4.23%  '
$style = $ws.Range("D21").Style
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '13.98'
$ws.Range("D21").Style = $style
$ws.Range("E21").Value = '  +2.46%  '
$style = $ws.Range("D22").Style
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.733'
$ws.Range("D22").Style = $style
$ws.Range("E22").Value = '  +2.18%  '
$style = $ws.Range("D23").Style
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '7.43'
$ws.Range("D23").Style = $style
$ws.Range("E23").Value = '  +4.74%  '
$style = $ws.Range("D24").Style
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '13.73'
$ws.Range("D24").Style = $style
$ws.Range("E24").Value = '  +3.33%  '
$style = $ws.Range("D25").Style
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '81.60'
$ws.Range("D25").Style = $style
$ws.Range("E25").Value = '  +1.04%  '
$style = $ws.Range("D26").Style
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.999'
$ws.Range("D26").Style = $style
$ws.Range("E26").Value = '  -0.13%  '
$ws.Range("E27").Value = '  +5.38%  '
$ws.Range("E28").Value = '  -0.05%  '
$ws.Range("E29").Value = '  +4.76%  '
$style = $ws.Range("D30").Style
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '8.19'
$ws.Range("D30").Style = $style
$ws.Range("E30").Value = '  +5.53%  '
$style = $ws.Range("D31").Style
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '6.77'
$ws.Range("D31").Style = $style
$ws.Range("E31").Value = '  +11.11%  '
$style = $ws.Range("D32").Style
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.114'
$ws.Range("D32").Style = $style
$ws.Range("E32").Value = '  +15.75%  '
$style = $ws.Range("D33").Style
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '26.78'
$ws.Range("D33").Style = $style
$ws.Range("E33").Value = '  +4.22%  '
$style = $ws.Range("D34").Style
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.04'
$ws.Range("D34").Style = $style
$ws.Range("E34").Value = '  +4.48%  '
$ws.Range("D35").Value = '0.0₃0790'
$ws.Range("E35").Value = '  +3.81%  '
$style = $ws.Range("D36").Style
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '6.04'
$ws.Range("D36").Style = $style
$ws.Range("E36").Value = '  +1.97%  '
$style = $ws.Range("D37").Style
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.19'
$ws.Range("D37").Style = $style
$ws.Range("E37").Value = '  +5.17%  '
$style = $ws.Range("D38").Style
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '50.13'
$ws.Range("D38").Style = $style
$ws.Range("E38").Value = '  +2.40%  '
$ws.Range("E39").Value = '  +8.81%  '
$ws.Range("E40").Value = '  +1.58%  '
$style = $ws.Range("D41").Style
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '420.51'
$ws.Range("D41").Style = $style
$ws.Range("E41").Value = '  +5.10%  '
$ws.Range("D42").Value = '2.950.05'
$ws.Range("E42").Value = '  +7.16%  '
$style = $ws.Range("D43").Style
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.0369'
$ws.Range("D43").Style = $style
$ws.Range("E43").Value = '  +5.23%  '
$ws.Range("E44").Value = '  +9.95%  '
$ws.Range("E45").Value = '  +1.02%  '
$style = $ws.Range("D46").Style
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.12'
$ws.Range("D46").Style = $style
$ws.Range("E46").Value = '  +6.03%  '
$style = $ws.Range("D48").Style
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '123.73'
$ws.Range("D48").Style = $style
$ws.Range("E48").Value = '  +2.00%  '
$style = $ws.Range("D49").Style
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '34.75'
$ws.Range("D49").Style = $style
$ws.Range("E49").Value = '  +0.12%  '
$ws.Range("E50").Value = '  +0.41%  '
$style = $ws.Range("D51").Style
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '24.33'
$ws.Range("D51").Style = $style
$ws.Range("E51").Value = '  +4.17%  '
